$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.327 (0.321 ± 0.005)"
$ws.Range("C2").Value = "00:04:47 (00:05:32 ± 00:00:48)"
$ws.Range("D2").Value = "00:00:01 (00:00:05 ± 00:00:02)"

$ws.Range("B3").Value = "0.379 (0.304 ± 0.039)"
$ws.Range("C3").Value = "00:00:33 (00:00:50 ± 00:00:21)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B4").Value = "0.351 (0.258 ± 0.037)"
$ws.Range("C4").Value = "00:00:26 (00:00:35 ± 00:00:08)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B5").Value = "0.349 (0.275 ± 0.043)"
$ws.Range("C5").Value = "00:05:05 (00:05:13 ± 00:00:04)"
$ws.Range("D5").Value = "00:00:01 (00:00:02 ± 00:00:01)"

$ws.Range("B6").Value = "0.363 (0.303 ± 0.032)"
$ws.Range("C6").Value = "00:04:58 (00:05:03 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 ± 00:00:00)"

$ws.Range("B9").Value = "0.369 (0.316 ± 0.040)"
$ws.Range("C9").Value = "00:04:59 (00:05:01 ± 00:00:02)"
$ws.Range("D9").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B11").Value = "0.281 (0.166 ± 0.070)"
$ws.Range("C11").Value = "00:05:05 (00:05:06 ± 00:00:00)"
$ws.Range("D11").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B12").Value = "0.294 (0.274 ± 0.016)"
$ws.Range("C12").Value = "00:01:12 (00:02:13 ± 00:00:36)"
$ws.Range("D12").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B13").Value = "0.094 (0.045 ± 0.023)"
$ws.Range("C13").Value = "00:00:02 (00:00:03 ± 00:00:01)"
$ws.Range("D13").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B14").Value = "0.344 (0.264 ± 0.043)"
$ws.Range("C14").Value = "00:01:48 (00:01:57 ± 00:00:04)"
$ws.Range("D14").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B15").Value = "0.358 (0.303 ± 0.032)"
$ws.Range("C15").Value = "00:00:46 (00:04:11 ± 00:01:23)"
$ws.Range("D15").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B16").Value = "0.389 (0.314 ± 0.043)"
$ws.Range("C16").Value = "00:10:02 (00:11:05 ± 00:00:53)"
$ws.Range("D16").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B17").Value = "0.385 (0.305 ± 0.033)"
$ws.Range("C17").Value = "00:05:01 (00:05:26 ± 00:00:17)"
$ws.Range("D17").Value = "00:00:00 (00:00:00 ± 00:00:00)"
